$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision resulted in no movie being selected for Friday.`n"
$ws.Range("D2").Value = "no_decision, "

$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be showcased on Friday.`n"
$ws.Range("D3").Value = "Barbie_was_selected, "

$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been confirmed.`n"
$ws.Range("D4").Value = "Oppenheimer_was_selected, "

$ws.Range("C5").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D5").Value = "both_movies, "

$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision about which movie will be shown on Friday was not reached, leading to no selection being made.`n"
$ws.Range("D6").Value = "no_decision, "

$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("D7").Value = "no_decision, "

$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("D8").Value = "Barbie_was_selected, "

$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision process concluded without an agreement on the movie to be shown on Friday, resulting in no decision being made.`n"
$ws.Range("D9").Value = "no_decision, "

$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday concluded with no final agreement on which film to show.`n"
$ws.Range("D10").Value = "no_decision, "

$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made.`n"
$ws.Range("D11").Value = "no_decision, "

$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("D12").Value = "both_movies, "

$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision about Friday's movie was not finalized, so no action will be taken regarding acquiring movie rights.`n"
$ws.Range("D13").Value = "no_decision, "

$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("D14").Value = "no_decision, "

$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie remains unresolved, so I have recorded that there was no decision made.`n"
$ws.Range("D15").Value = "no_decision, "

$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded, reflecting that no movie will be shown on Friday, as the committee did not reach a consensus.`n"
$ws.Range("D16").Value = "no_decision, "

$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("D17").Value = "Barbie_was_selected, "

$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D18").Value = "no_decision, "

$ws.Range("C19").Value = "MSG: None`n`nMSG: I have successfully acquired the rights to both movies.`n"
$ws.Range("D19").Value = "both_movies, "

$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("D20").Value = "both_movies, "

$ws.Range("C21").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding a movie for Friday.`n"
$ws.Range("D21").Value = "no_decision, "

$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded and `"Barbie`" will be the movie acquired for Friday.`n"
$ws.Range("D22").Value = "Barbie_was_selected, "

$ws.Range("C23").Value = "MSG: None`n`nMSG: The function for no decision has been called, indicating that there was no agreement reached on which movie to show on Friday.`n"
$ws.Range("D23").Value = "no_decision, "

$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been made: there is no definitive choice for the movie to be shown on Friday.`n"
$ws.Range("D24").Value = "no_decision, "

$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has been recorded as no decision being made.`n"
$ws.Range("D25").Value = "no_decision, "

$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been finalized.`n"
$ws.Range("D26").Value = "Barbie_was_selected, "

$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has resulted in no agreement. Therefore, no further action is taken.`n"
$ws.Range("D27").Value = "no_decision, "

$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded; the movie `"Barbie`" will be acquired for the event.`n"
$ws.Range("D28").Value = "Barbie_was_selected, "

$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D29").Value = "no_decision, "

$ws.Range("C30").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for Friday's showing.`n"
$ws.Range("D30").Value = "Barbie_was_selected, "

$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been finalized.`n"
$ws.Range("D31").Value = "Barbie_was_selected, "

$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision-making process concluded without a consensus, resulting in no movie being selected for Friday.`n"
$ws.Range("D32").Value = "no_decision, "

$ws.Range("C33").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("D33").Value = "no_decision, "

$ws.Range("C34").Value = "MSG: None`n`nMSG: The committee did not reach a decision on what movie to show on Friday, so the selection process is concluded without a choice.`n"
$ws.Range("D34").Value = "no_decision, "

$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision to select a movie for Friday resulted in no final agreement, so the outcome is marked as `"no decision.`"`n"
$ws.Range("D35").Value = "no_decision, "

$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday remains unresolved.`n"
$ws.Range("D36").Value = "no_decision, "

$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("D37").Value = "Barbie_was_selected, "

$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been made.`n"
$ws.Range("D38").Value = "Barbie_was_selected, "

$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision process concluded without a choice of movie for Friday, indicating that no decision was made.`n"
$ws.Range("D39").Value = "no_decision, "

$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be selected for Friday.`n"
$ws.Range("D40").Value = "no_decision, "

$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded, and it has been determined that no movie was selected for Friday.`n"
$ws.Range("D41").Value = "no_decision, "
